# Apply updated per-unit voltage magnitude results for the 380 kV case
# (commit: "case with 380 kV done") to Code/Results/Cases/Case_2_104/res_bus/vm_pu.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.030976068851771
$ws.Range("D2").Value = 1.034376972943943
$ws.Range("E2").Value = 1.0511234155851
$ws.Range("F2").Value = 1.056242480241024
$ws.Range("I2").Value = 1.036280307663609
$ws.Range("J2").Value = 1.03611434355409
$ws.Range("K2").Value = 1.037176622660992
$ws.Range("L2").Value = 1.053875790400496
$ws.Range("M2").Value = 1.058980728218371
$ws.Range("N2").Value = 1.037585744765914

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.031807927830352
$ws.Range("D3").Value = 1.034984485751495
$ws.Range("E3").Value = 1.052342410517915
$ws.Range("F3").Value = 1.057521963342256
$ws.Range("I3").Value = 1.036474787651092
$ws.Range("J3").Value = 1.036588482674507
$ws.Range("K3").Value = 1.03759384806658
$ws.Range("L3").Value = 1.054906323094512
$ws.Range("M3").Value = 1.060072625578378
$ws.Range("N3").Value = 1.038060557218266

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.032346186959689
$ws.Range("D4").Value = 1.035377485179999
$ws.Range("E4").Value = 1.053132044903663
$ws.Range("F4").Value = 1.058350709269872
$ws.Range("I4").Value = 1.036599200174632
$ws.Range("J4").Value = 1.036894617177478
$ws.Range("K4").Value = 1.037863025174119
$ws.Range("L4").Value = 1.055573442731532
$ws.Range("M4").Value = 1.060779435026419
$ws.Range("N4").Value = 1.038367126467345

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.032572467857674
$ws.Range("D5").Value = 1.035542676419019
$ws.Range("E5").Value = 1.053464214704023
$ws.Range("F5").Value = 1.058699314770368
$ws.Range("I5").Value = 1.036651160531115
$ws.Range("J5").Value = 1.037023156099398
$ws.Range("K5").Value = 1.037975995906458
$ws.Range("L5").Value = 1.055853971423884
$ws.Range("M5").Value = 1.06107664512071
$ws.Range("N5").Value = 1.038495847929277

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.032610461196336
$ws.Range("D6").Value = 1.03557041120857
$ws.Range("E6").Value = 1.053519999668574
$ws.Range("F6").Value = 1.058757858962843
$ws.Range("I6").Value = 1.036659864793208
$ws.Range("J6").Value = 1.037044728943731
$ws.Range("K6").Value = 1.037994952945477
$ws.Range("L6").Value = 1.055901077645158
$ws.Range("M6").Value = 1.061126551974775
$ws.Range("N6").Value = 1.038517451409524

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.032349210550139
$ws.Range("D7").Value = 1.035379692576136
$ws.Range("E7").Value = 1.053136482556346
$ws.Range("F7").Value = 1.05835536656388
$ws.Range("I7").Value = 1.036599895818607
$ws.Range("J7").Value = 1.036896335350596
$ws.Range("K7").Value = 1.037864535447416
$ws.Range("L7").Value = 1.055577190887415
$ws.Range("M7").Value = 1.060783406096582
$ws.Range("N7").Value = 1.038368847080466

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.031257200789263
$ws.Range("D8").Value = 1.034582305073859
$ws.Range("E8").Value = 1.051535202426765
$ws.Range("F8").Value = 1.056674715208688
$ws.Range("I8").Value = 1.036346328643361
$ws.Range("J8").Value = 1.036274718364392
$ws.Range("K8").Value = 1.037317790171353
$ws.Range("L8").Value = 1.0542240028746
$ws.Range("M8").Value = 1.059349682694803
$ws.Range("N8").Value = 1.03774634732684

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.029332914146692
$ws.Range("D9").Value = 1.033176482754423
$ws.Range("E9").Value = 1.048720114238958
$ws.Range("F9").Value = 1.053719547881194
$ws.Range("I9").Value = 1.035888590563603
$ws.Range("J9").Value = 1.035174288829753
$ws.Range("K9").Value = 1.036348298802715
$ws.Range("L9").Value = 1.051841748162252
$ws.Range("M9").Value = 1.056825385407449
$ws.Range("N9").Value = 1.036644355056041

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.028050098020781
$ws.Range("D10").Value = 1.032238857854666
$ws.Range("E10").Value = 1.046847751065781
$ws.Range("F10").Value = 1.051753652873573
$ws.Range("I10").Value = 1.035576121840603
$ws.Range("J10").Value = 1.034437310868449
$ws.Range("K10").Value = 1.035697948590033
$ws.Range("L10").Value = 1.050255034417607
$ws.Range("M10").Value = 1.055143893895154
$ws.Range("N10").Value = 1.0359063305015

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.027494646306712
$ws.Range("D11").Value = 1.031832774544587
$ws.Range("E11").Value = 1.046038018573909
$ws.Range("F11").Value = 1.050903385451417
$ws.Range("I11").Value = 1.035439092166161
$ws.Range("J11").Value = 1.034117403765293
$ws.Range("K11").Value = 1.035415395759557
$ws.Range("L11").Value = 1.049568304112455
$ws.Range("M11").Value = 1.054416104944518
$ws.Range("N11").Value = 1.035585969093564

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.027288330134949
$ws.Range("D12").Value = 1.031681925520368
$ws.Range("E12").Value = 1.045737398856311
$ws.Range("F12").Value = 1.050587703981445
$ws.Range("I12").Value = 1.035387934020217
$ws.Range("J12").Value = 1.033998457896956
$ws.Range("K12").Value = 1.035310301569679
$ws.Range("L12").Value = 1.049313270021594
$ws.Range("M12").Value = 1.054145817011622
$ws.Range("N12").Value = 1.03546685430845

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.027332585489154
$ws.Range("D13").Value = 1.031714283653145
$ws.Range("E13").Value = 1.045801875993545
$ws.Range("F13").Value = 1.050655412154559
$ws.Range("I13").Value = 1.035398919338783
$ws.Range("J13").Value = 1.034023977524063
$ws.Range("K13").Value = 1.035332851020416
$ws.Range("L13").Value = 1.049367973538877
$ws.Range("M13").Value = 1.054203792657059
$ws.Range("N13").Value = 1.035492410176355

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.027477592074523
$ws.Range("D14").Value = 1.031820305546659
$ws.Range("E14").Value = 1.046013166181335
$ws.Range("F14").Value = 1.050877288157359
$ws.Range("I14").Value = 1.035434868703125
$ws.Range("J14").Value = 1.034107574068529
$ws.Range("K14").Value = 1.035406711522892
$ws.Range("L14").Value = 1.049547221939283
$ws.Range("M14").Value = 1.054393761941459
$ws.Range("N14").Value = 1.035576125437503

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.027566935885333
$ws.Range("D15").Value = 1.031885627634315
$ws.Range("E15").Value = 1.046143369009643
$ws.Range("F15").Value = 1.051014012556402
$ws.Range("I15").Value = 1.035456983957894
$ws.Range("J15").Value = 1.034159065026015
$ws.Range("K15").Value = 1.035452200690698
$ws.Range("L15").Value = 1.04965766914023
$ws.Range("M15").Value = 1.054510814259223
$ws.Range("N15").Value = 1.035627689518054

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.028086961938688
$ws.Range("D16").Value = 1.032265806571098
$ws.Range("E16").Value = 1.046901511584917
$ws.Range("F16").Value = 1.051810102823513
$ws.Range("I16").Value = 1.035585179669636
$ws.Range("J16").Value = 1.034458525446496
$ws.Range("K16").Value = 1.035716680806816
$ws.Range("L16").Value = 1.050300617251513
$ws.Range("M16").Value = 1.05519220125564
$ws.Range("N16").Value = 1.035927575206681

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.028413165500138
$ws.Range("D17").Value = 1.032504260756226
$ws.Range("E17").Value = 1.047377344766059
$ws.Range("F17").Value = 1.052309730202687
$ws.Range("I17").Value = 1.035665130986969
$ws.Range("J17").Value = 1.034646157822393
$ws.Range("K17").Value = 1.035882329231474
$ws.Range("L17").Value = 1.050704008485521
$ws.Range("M17").Value = 1.05561969907533
$ws.Range("N17").Value = 1.036115474042074

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.028603436014863
$ws.Range("D18").Value = 1.032643338711145
$ws.Range("E18").Value = 1.047654988043323
$ws.Range("F18").Value = 1.052601248854115
$ws.Range("I18").Value = 1.035711598414898
$ws.Range("J18").Value = 1.034755524268265
$ws.Range("K18").Value = 1.035978857741595
$ws.Range("L18").Value = 1.050939331408633
$ws.Range("M18").Value = 1.055869081101071
$ws.Range("N18").Value = 1.036224995800843

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.028668313569663
$ws.Range("D19").Value = 1.032690759247478
$ws.Range("E19").Value = 1.0477496738575
$ws.Range("F19").Value = 1.052700665239537
$ws.Range("I19").Value = 1.035727414296114
$ws.Range("J19").Value = 1.034792802429954
$ws.Range("K19").Value = 1.036011755940169
$ws.Range("L19").Value = 1.051019575842706
$ws.Range("M19").Value = 1.055954119008388
$ws.Range("N19").Value = 1.036262326901798

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.028378166783812
$ws.Range("D20").Value = 1.032478677732136
$ws.Range("E20").Value = 1.047326282258292
$ws.Range("F20").Value = 1.052256115152236
$ws.Range("I20").Value = 1.03565657021899
$ws.Range("J20").Value = 1.03462603453591
$ws.Range("K20").Value = 1.035864566172892
$ws.Range("L20").Value = 1.050660725158321
$ws.Range("M20").Value = 1.055573829531685
$ws.Range("N20").Value = 1.036095322178216

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.027434891150264
$ws.Range("D21").Value = 1.031789085051812
$ws.Range("E21").Value = 1.045950942329247
$ws.Range("F21").Value = 1.05081194719229
$ws.Range("I21").Value = 1.035424289663273
$ws.Range("J21").Value = 1.034082960225192
$ws.Range("K21").Value = 1.03538496534541
$ws.Range("L21").Value = 1.049494436432957
$ws.Range("M21").Value = 1.054337819479742
$ws.Range("N21").Value = 1.035551476639686

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.026841836011818
$ws.Range("D22").Value = 1.031355444252566
$ws.Range("E22").Value = 1.045087084467782
$ws.Range("F22").Value = 1.04990478389387
$ws.Range("I22").Value = 1.035276746111406
$ws.Range("J22").Value = 1.033740824743607
$ws.Range("K22").Value = 1.035082603044288
$ws.Range("L22").Value = 1.048761422324495
$ws.Range("M22").Value = 1.053560952444252
$ws.Range("N22").Value = 1.035208855286473

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.027156223594275
$ws.Range("D23").Value = 1.031585331333814
$ws.Range("E23").Value = 1.04554494952862
$ws.Range("F23").Value = 1.050385608883285
$ws.Range("I23").Value = 1.035355103717165
$ws.Range("J23").Value = 1.0339222617602
$ws.Range("K23").Value = 1.035242968322193
$ws.Range("L23").Value = 1.049149980976428
$ws.Range("M23").Value = 1.05397275989391
$ws.Range("N23").Value = 1.035390549964441

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.028393981192636
$ws.Range("D24").Value = 1.032490237625393
$ws.Range("E24").Value = 1.047349354904934
$ws.Range("F24").Value = 1.05228034119347
$ws.Range("I24").Value = 1.035660438977027
$ws.Range("J24").Value = 1.034635127618799
$ws.Range("K24").Value = 1.035872592816968
$ws.Range("L24").Value = 1.050680282932682
$ws.Range("M24").Value = 1.055594555911684
$ws.Range("N24").Value = 1.036104428174326

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.029830385347185
$ws.Range("D25").Value = 1.03353999937109
$ws.Range("E25").Value = 1.049447109020472
$ws.Range("F25").Value = 1.054482782172385
$ws.Range("I25").Value = 1.036008217620251
$ws.Range("J25").Value = 1.035459371404121
$ws.Range("K25").Value = 1.036599648174116
$ws.Range("L25").Value = 1.057477731570414
$ws.Range("M25").Value = 1.058980728218371
$ws.Range("N25").Value = 1.036929842480365
